$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# Push the existing entries down one row (row 8's data -> row 9, ... row 3's data -> row 4)
# so the newest result can be written into row 3. Walk bottom-up so a source row is always
# read before it gets overwritten.
for ($r = 8; $r -ge 3; $r--) {
    $src = $r
    $dst = $r + 1
    $ws.Range("A$dst").Value = $ws.Range("A$src").Value()
    $ws.Range("B$dst").Value = $ws.Range("B$src").Value()
    $ws.Range("C$dst").Value = $ws.Range("C$src").Value()
}

# Write the newest 4D box result into the now-vacated row 3
$ws.Range("A3").Value = "28/6/2025 (Sat)"
$ws.Range("B3").Value = "3 4 6 0" + $nl + "4 9 3 6" + $nl + "1 5 2 7" + $nl + "0 0 4 8"
$ws.Range("C3").Value = "✅ Direct: 14/3980 (0.35%)" + $nl + "✅ iBet: 14/215 (6.51%)"

# Row 12 picks up a matching (still-empty) placeholder cell in column C
$ws.Cells.Item(12, 3).WrapText = $true

# A fresh blank placeholder row is appended at the bottom of the table
$ws.Cells.Item(38, 2).WrapText = $true
